# Sorption isotherm results rewritten
# Adds five new data rows (20-24) for the "Sodium Montmorillonite" /
# "Tamamura" source to the Compilation sheet's results table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compilation")

# --- Row 20 -----------------------------------------------------------
$ws.Range("A20").Value   = "Sodium Montmorillonite"
$ws.Range("B20").Formula = "=0.03/100"
$ws.Range("C20").Value   = "10 mM NaCl"
$ws.Range("D20").Value   = 3
$ws.Range("E20").Value   = 6740.15
$ws.Range("F20").Value   = "Experimental"

# --- Rows 21-23 (share the B-column formula, like the rows above them) -
$ws.Range("B21:B23").Formula = "=0.03/100"

$ws.Range("A21").Value = "Sodium Montmorillonite"
$ws.Range("C21").Value = "10 mM NaCl"
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 17749.39
$ws.Range("F21").Value = "Experimental"

$ws.Range("A22").Value = "Sodium Montmorillonite"
$ws.Range("C22").Value = "10 mM NaCl"
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = 21473.27
$ws.Range("F22").Value = "Experimental"

$ws.Range("A23").Value = "Sodium Montmorillonite"
$ws.Range("C23").Value = "10 mM NaCl"
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = 22894.86
$ws.Range("F23").Value = "Experimental"

# --- Row 24 (new source "Tamamura") ------------------------------------
$ws.Range("A24").Value   = "Sodium Montmorillonite"
$ws.Range("B24").Formula = "=0.1/30"
$ws.Range("C24").Value   = "10 mM NaCl"
$ws.Range("D24").Value   = 5.25
$ws.Range("E24").Formula = "=1/(0.03/30*0.1/0.97)"
$ws.Range("F24").Value   = "Tamamura"

# Leave the same cell selected as in the authored workbook.
[void]$ws.Range("B11").Select()
